$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: 100n capacitor - quantity and reference designators change ---
$ws.Range("B2").Value = 2
$ws.Range("G2").Value = "C3, C4, C5"

# --- Row 3: 2n2 capacitor - quantity and reference designators change ---
$ws.Range("B3").Value = 3
$ws.Range("G3").Value = "C1, C2"

# --- Row 10: was the MH18-1 header connector, becomes the MagJack connector ---
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "11A009-1702-00"
$ws.Range("D10").Value = "MENTECH"
$ws.Range("E10").Value = "11A009-1702-00"
$ws.Range("F10").Value = "11A009-1702-00"
$ws.Range("G10").Value = "J1"
$ws.Range("H10").Value = "CONN MAGJACK 1PORT 100 BASE-T"

# --- Row 11: was the MagJack connector, becomes the MH18-1 header connector ---
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "MH18-1"
$ws.Range("D11").Value = "ANY"
$ws.Range("E11").Value = "MH18-1-0.1"
$ws.Range("F11").Value = "MH18-1-0.1"
$ws.Range("G11").Value = "J2, J3"
$ws.Range("H11").Value = "CONN HEADER VERT .100 1ROW 18POS 8.08 HEAD 3.05 TAIL 15AU"

# --- Row 13: was 150R/R12, becomes 330R/R11,R12,R13 ---
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "330R"
$ws.Range("D13").Value = "ANY"
$ws.Range("E13").Value = "R0603_330R_5%_125mW"
$ws.Range("F13").Value = "R0603"
$ws.Range("G13").Value = "R11, R12, R13"
$ws.Range("H13").Value = "RES 330 OHM 1/8W 5% 0603 SMD"

# --- Row 14: was 330R/R5,R13,R14, becomes 270R/R5..R10 ---
$ws.Range("B14").Value = 6
$ws.Range("C14").Value = "270R"
$ws.Range("D14").Value = "ANY"
$ws.Range("E14").Value = "R0603_270R_5%_125mW"
$ws.Range("F14").Value = "R0603"
$ws.Range("G14").Value = "R5, R6, R7, R8, R9, R10"
$ws.Range("H14").Value = "RES 270 OHM 1/8W 5% 0603 SMD"

# --- Row 15 (old trailing 270R row) is removed entirely; table now ends at row 14 ---
$ws.Rows("15").Delete()
